$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 26)
$ws.Range("D2").Value = 0.9999988332380492
$ws.Range("E2").Value = 0.9999988332380492

# Row 3 (Control 33)
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.9999999977660445
$ws.Range("E3").Value = 0.9999999977660445

# Row 4 (Control 36)
$ws.Range("D4").Value = 0.9804776950544825
$ws.Range("E4").Value = 0.9804776950544825

# Row 5 (Control 49)
$ws.Range("D5").Value = 0.001058286236534054
$ws.Range("E5").Value = 0.001058286236534054

# Row 6 (Control 2)
$ws.Range("D6").Value = [double]"1.68652121484057E-12"
$ws.Range("E6").Value = [double]"1.68652121484057E-12"

# Row 7 (MDD 42)
$ws.Range("D7").Value = [double]"9.699983565650095E-09"
$ws.Range("E7").Value = 0.9999999903000164

# Row 8 (MDD 52)
$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = [double]"2.220446049250313E-16"

# Row 9 (MDD 20)
$ws.Range("D9").Value = 0.9930961941136397
$ws.Range("E9").Value = 0.006903805886360304

# Row 10 (MDD 51)
$ws.Range("D10").Value = [double]"8.884696294575854E-06"
$ws.Range("E10").Value = 0.9999911153037054

# Row 11 (MDD 40)
$ws.Range("D11").Value = [double]"1.02410861889238E-09"
$ws.Range("E11").Value = 0.9999999989758914
$ws.Range("F11").Value = 8.830671310424805
$ws.Range("G11").Value = 0.4
